$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text storage for
# numeric-looking strings (prices, percentages, etc.) so that
# "306.40" does not get silently coerced into the float 306.4 and
# lose its trailing zero / sign / percent-suffix. A leading apostrophe
# forces Excel to keep the literal text; re-applying the "Normal" style
# afterwards clears the transient quote-prefix formatting flag that
# operation leaves behind, so the stored cell style stays unchanged.
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '306.40'
Set-TextValue $ws.Range("E2") '-0.14%'
Set-TextValue $ws.Range("D3") '40.30'
Set-TextValue $ws.Range("E3") '0.97%'
Set-TextValue $ws.Range("D4") '5.115'
Set-TextValue $ws.Range("E4") '0.94%'
Set-TextValue $ws.Range("D5") '0.07586'
Set-TextValue $ws.Range("E5") '-2.35%'
Set-TextValue $ws.Range("D6") '1.610'
Set-TextValue $ws.Range("E6") '-2.12%'
Set-TextValue $ws.Range("B7") 'MXToken'
Set-TextValue $ws.Range("C7") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D7") '0.9064'
Set-TextValue $ws.Range("E7") '-1.28%'
Set-TextValue $ws.Range("B8") 'BTSEToken'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D8") '2.423'
Set-TextValue $ws.Range("E8") '-5.35%'
Set-TextValue $ws.Range("D9") '0.1010'
Set-TextValue $ws.Range("E9") '3.80%'
Set-TextValue $ws.Range("D10") '0.1749'
Set-TextValue $ws.Range("E10") '0.61%'
Set-TextValue $ws.Range("D11") '0.09089'
Set-TextValue $ws.Range("E11") '1.70%'
Set-TextValue $ws.Range("D12") '0.04211'
Set-TextValue $ws.Range("E12") '-4.16%'
Set-TextValue $ws.Range("D13") '0.1055'
Set-TextValue $ws.Range("E13") '-0.37%'
Set-TextValue $ws.Range("D14") '0.001245'
Set-TextValue $ws.Range("E14") '-2.36%'
Set-TextValue $ws.Range("D15") '0.005875'
Set-TextValue $ws.Range("E15") '3.78%'
Set-TextValue $ws.Range("D16") '3.354'
Set-TextValue $ws.Range("E16") '-0.39%'
Set-TextValue $ws.Range("D17") '4.271'
Set-TextValue $ws.Range("E17") '-1.14%'
Set-TextValue $ws.Range("E18") '-2.79%'
Set-TextValue $ws.Range("D19") '6.643'
Set-TextValue $ws.Range("E19") '-5.66%'
Set-TextValue $ws.Range("D20") '0.1357'
Set-TextValue $ws.Range("E20") '-0.54%'
Set-TextValue $ws.Range("E21") '2.50%'
Set-TextValue $ws.Range("D22") '0.04178'
Set-TextValue $ws.Range("E22") '0.65%'
Set-TextValue $ws.Range("D23") '0.001227'
Set-TextValue $ws.Range("E23") '1.67%'
Set-TextValue $ws.Range("D24") '0.004054'
Set-TextValue $ws.Range("E24") '-0.78%'
Set-TextValue $ws.Range("E25") '6.07%'
Set-TextValue $ws.Range("D26") '0.0003012'
Set-TextValue $ws.Range("E26") '0.49%'
Set-TextValue $ws.Range("D38") '0.02385'
Set-TextValue $ws.Range("E38") '0.39%'
Set-TextValue $ws.Range("D39") '0.05138'
Set-TextValue $ws.Range("E39") '-0.41%'
Set-TextValue $ws.Range("D40") '0.007781'
Set-TextValue $ws.Range("E40") '-2.43%'
Set-TextValue $ws.Range("D41") '0.1298'
Set-TextValue $ws.Range("E41") '-2.22%'
Set-TextValue $ws.Range("D42") '0.007057'
Set-TextValue $ws.Range("E42") '-6.70%'
Set-TextValue $ws.Range("D43") '0.001921'
Set-TextValue $ws.Range("E43") '-4.91%'
Set-TextValue $ws.Range("D44") '0.008452'
Set-TextValue $ws.Range("E44") '4.78%'
Set-TextValue $ws.Range("D45") '0.3308'
Set-TextValue $ws.Range("E45") '-0.57%'
Set-TextValue $ws.Range("D46") '0.00006368'
Set-TextValue $ws.Range("E46") '-5.57%'
Set-TextValue $ws.Range("D47") '0.00000000751'
Set-TextValue $ws.Range("E47") '-0.44%'
Set-TextValue $ws.Range("B48") 'BOLO'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws.Range("D48") '0.006845'
Set-TextValue $ws.Range("E48") '99.75%'
Set-TextValue $ws.Range("B49") 'CoinbaseStockToken'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws.Range("D49") '0.004407'
Set-TextValue $ws.Range("E49") '6.83%'
Set-TextValue $ws.Range("D50") '0.00002103'
Set-TextValue $ws.Range("E50") '-0.44%'
Set-TextValue $ws.Range("D51") '0.0002003'
Set-TextValue $ws.Range("E51") '-0.44%'
